$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '26.674.79'
Set-TextValue $ws.Range("E2") '  -2.65%  '
Set-TextValue $ws.Range("D3") '1.554.27'
Set-TextValue $ws.Range("E3") '  -0.59%  '
Set-TextValue $ws.Range("E4") '  +0.13%  '
Set-TextValue $ws.Range("D5") '205.49'
Set-TextValue $ws.Range("E5") '  -1.16%  '
Set-TextValue $ws.Range("E6") '  -1.97%  '
Set-TextValue $ws.Range("E7") '  +0.09%  '
Set-TextValue $ws.Range("D8") '21.77'
Set-TextValue $ws.Range("E8") '  -0.18%  '
Set-TextValue $ws.Range("E9") '  -0.67%  '
Set-TextValue $ws.Range("D10") '0.0581'
Set-TextValue $ws.Range("E10") '  -1.35%  '
Set-TextValue $ws.Range("D11") '0.0863'
Set-TextValue $ws.Range("E11") '  -0.58%  '
Set-TextValue $ws.Range("D12") '1.778.65'
Set-TextValue $ws.Range("E12") '  -0.44%  '
Set-TextValue $ws.Range("D13") '1.556.42'
Set-TextValue $ws.Range("E13") '  -1.10%  '
Set-TextValue $ws.Range("D14") '3.72'
Set-TextValue $ws.Range("E14") '  -2.51%  '
Set-TextValue $ws.Range("D15") '0.510'
Set-TextValue $ws.Range("E15") '  -0.46%  '
Set-TextValue $ws.Range("D16") '61.46'
Set-TextValue $ws.Range("E16") '  -2.86%  '
Set-TextValue $ws.Range("D17") '26.725.77'
Set-TextValue $ws.Range("E17") '  -2.48%  '
Set-TextValue $ws.Range("D18") '213.44'
Set-TextValue $ws.Range("E18") '  +0.75%  '
Set-TextValue $ws.Range("E19") '  +0.95%  '
Set-TextValue $ws.Range("D20") '0.0₃0674'
Set-TextValue $ws.Range("E20") '  -1.91%  '
Set-TextValue $ws.Range("E21") '  +0.08%  '
Set-TextValue $ws.Range("D22") '4.07'
Set-TextValue $ws.Range("E22") '  -0.88%  '
Set-TextValue $ws.Range("D23") '9.36'
Set-TextValue $ws.Range("E23") '  -1.62%  '
Set-TextValue $ws.Range("E24") '  +0.16%  '
Set-TextValue $ws.Range("D25") '152.73'
Set-TextValue $ws.Range("E25") '  -0.31%  '
Set-TextValue $ws.Range("D26") '6.73'
Set-TextValue $ws.Range("E26") '  +0.68%  '
Set-TextValue $ws.Range("D27") '14.81'
Set-TextValue $ws.Range("E27") '  -0.97%  '
Set-TextValue $ws.Range("E28") '  +0.13%  '
Set-TextValue $ws.Range("E29") '  -1.09%  '
Set-TextValue $ws.Range("E30") '  -1.47%  '
Set-TextValue $ws.Range("D31") '1.11'
Set-TextValue $ws.Range("E31") '  -4.05%  '
Set-TextValue $ws.Range("D32") '3.14'
Set-TextValue $ws.Range("E32") '  -1.46%  '
Set-TextValue $ws.Range("D33") '1.377.89'
Set-TextValue $ws.Range("E33") '  +1.17%  '
Set-TextValue $ws.Range("E34") '  -1.40%  '
Set-TextValue $ws.Range("D35") '1.55'
Set-TextValue $ws.Range("E35") '  +0.98%  '
Set-TextValue $ws.Range("E36") '  -0.56%  '
Set-TextValue $ws.Range("D37") '0.935'
Set-TextValue $ws.Range("E37") '  -3.87%  '
Set-TextValue $ws.Range("E38") '  -2.57%  '
Set-TextValue $ws.Range("D39") '0.517'
Set-TextValue $ws.Range("E39") '  -2.45%  '
Set-TextValue $ws.Range("D40") '0.807'
Set-TextValue $ws.Range("E40") '  -1.52%  '
Set-TextValue $ws.Range("E41") '  +0.08%  '
Set-TextValue $ws.Range("D42") '0.993'
Set-TextValue $ws.Range("E42") '  +2.04%  '
Set-TextValue $ws.Range("E43") '  +2.21%  '
Set-TextValue $ws.Range("D46") '62.94'
Set-TextValue $ws.Range("E46") '  -1.65%  '
Set-TextValue $ws.Range("D47") '1.690.80'
Set-TextValue $ws.Range("E47") '  -0.51%  '
Set-TextValue $ws.Range("D48") '85.42'
Set-TextValue $ws.Range("E48") '  -0.08%  '
Set-TextValue $ws.Range("E49") '  -2.33%  '
Set-TextValue $ws.Range("E50") '  -0.05%  '
Set-TextValue $ws.Range("D51") '0.0942'
Set-TextValue $ws.Range("E51") '  -1.13%  '

# Row 44/45 swap (RenderToken <-> MXToken) with updated data
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D44") '2.17'
Set-TextValue $ws.Range("E44") '  +1.30%  '

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D45") '1.76'
Set-TextValue $ws.Range("E45") '  -1.62%  '
